$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F (used for the "수집일" / collected-date timestamps)
# so the full "YYYY-MM-DD HH:MM:SS" value is visible.
$ws.Columns.Item(6).ColumnWidth = 35.8

# Bump every F2:F73 timestamp forward by exactly one day
# (45677.56619212963 -> 45678.56619212963), keeping the same time-of-day.
for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value2 = $cell.Value2 + 1
}
